# Apply weekly update: insert two new price records for "Pepino ensalada"
# (Vega Modelo de Temuco) right after the existing row 824, shifting all
# subsequent rows down by two. This matches the diff where the dataset
# dimension grows from A1:R902 to A1:R904.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 825 (rows below shift down)
$ws.Range("A825:A826").EntireRow.Insert()

# ---- Row 825 : new record ----
$ws.Range("A825").Value = 10
$ws.Range("B825").Value = "Vega Modelo de Temuco"
$ws.Range("C825").Value = "La Araucanía"
$ws.Range("D825").Value = 45194
$ws.Range("E825").Value = 9
$ws.Range("F825").Value = 100112043
$ws.Range("G825").Value = "Pepino ensalada"
$ws.Range("H825").Value = "Sin especificar"
$ws.Range("I825").Value = "Primera"
$ws.Range("J825").Value = 650
$ws.Range("K825").Value = 13000
$ws.Range("L825").Value = 14000
$ws.Range("M825").Value = 13462
$ws.Range("N825").Value = "`$/caja 50 unidades"
$ws.Range("O825").Value = "Región de Arica y Parinacota"
$ws.Range("P825").Value = 269
$ws.Range("Q825").Value = 50
$ws.Range("R825").Value = "Hortaliza"

# ---- Row 826 : new record ----
$ws.Range("A826").Value = 10
$ws.Range("B826").Value = "Vega Modelo de Temuco"
$ws.Range("C826").Value = "La Araucanía"
$ws.Range("D826").Value = 45194
$ws.Range("E826").Value = 9
$ws.Range("F826").Value = 100112043
$ws.Range("G826").Value = "Pepino ensalada"
$ws.Range("H826").Value = "Sin especificar"
$ws.Range("I826").Value = "Segunda"
$ws.Range("J826").Value = 300
$ws.Range("K826").Value = 11000
$ws.Range("L826").Value = 12000
$ws.Range("M826").Value = 11667
$ws.Range("N826").Value = "`$/caja 60 unidades"
$ws.Range("O826").Value = "Región de Arica y Parinacota"
$ws.Range("P826").Value = 194
$ws.Range("Q826").Value = 60
$ws.Range("R826").Value = "Hortaliza"
